# Apply the "improving data, I hope" edit:
#  - delete the empty "2014_Nov" sheet
#  - rename "2015_Apr" -> "Selected"
#  - fix a typo in J4's note (missing colon)
#  - narrow columns E/F slightly
#  - add four new award rows (14-17) with the Nov/Apr 2014 data,
#    bolding the new Organisation names (F14:F16)
#  - move the saved selection to I19 (matches the diff's <selection>)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2015_Apr")

# --- fix existing cell J4: add missing colon after "available" ---
$ws.Range("J4").Value = "2013 application available: http://www.mhs.manchester.ac.uk/media/mhs/mhswebteam/documents/athena-swan/bronze2013.pdf"

# --- narrow columns E and F (stored widths end up 5.5 / 56 after Excel's padding) ---
$ws.Columns.Item(5).ColumnWidth = 4.666666666666667
$ws.Columns.Item(6).ColumnWidth = 55.166666666666664

# --- new row 14: Plymouth University (New, Bronze, Nov 2014) ---
$ws.Range("A14").Value = "Dept"
$ws.Range("B14").Value = "New"
$ws.Range("C14").Value = "Bronze"
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = "Nov"
$ws.Range("F14").Value = "Plymouth University – Peninsula Schools of Medicine and Dentistry"
$ws.Range("F14").Font.Bold = $true
$ws.Range("F14").Font.Size = 11
$ws.Range("G14").Value = "https://www.plymouth.ac.uk/your-university/about-us/university-structure/faculties/medicine-dentistry/athena-swan-at-pu-psmd"
$ws.Range("H14").Value = "Yes"
$ws.Range("I14").Value = "https://www.plymouth.ac.uk/uploads/production/document/path/4/4907/PU_PSMD_AS_Bronze_application_November_2014_web_version.pdf"

# --- new row 15: University of Leeds (Renewal, Bronze, Nov 2014) ---
$ws.Range("A15").Value = "Dept"
$ws.Range("B15").Value = "Renewal"
$ws.Range("C15").Value = "Bronze"
$ws.Range("D15").Value = 2014
$ws.Range("E15").Value = "Nov"
$ws.Range("F15").Value = "University of Leeds – School of Medicine "
$ws.Range("F15").Font.Bold = $true
$ws.Range("F15").Font.Size = 11
$ws.Range("G15").Value = "http://medhealth.leeds.ac.uk/homepage/516/school_of_medicine-athena_swan"
$ws.Range("H15").Value = "Yes"
$ws.Range("I15").Value = "http://medhealth.leeds.ac.uk/info/302/athena_swan/2015/applications"

# --- new row 16: Imperial College London (New, Silver, Nov 2014) ---
$ws.Range("A16").Value = "Dept"
$ws.Range("B16").Value = "New"
$ws.Range("C16").Value = "Silver"
$ws.Range("D16").Value = 2014
$ws.Range("E16").Value = "Nov"
$ws.Range("F16").Value = "Imperial College London – Department of Medicine"
$ws.Range("F16").Font.Bold = $true
$ws.Range("F16").Font.Size = 11
$ws.Range("G16").Value = "https://www1.imperial.ac.uk/departmentofmedicine/new_dom_life/dom_athena_swan/"
$ws.Range("H16").Value = "Yes"
$ws.Range("I16").Value = "https://www1.imperial.ac.uk/resources/FB2FDB0B-5F2A-48D0-93C8-829254CA6824/departmentofmedicineimperialsilver1.pdf"

# --- new row 17: University of Nottingham (New, Bronze, Apr 2014) ---
$ws.Range("A17").Value = "Dept"
$ws.Range("B17").Value = "New"
$ws.Range("C17").Value = "Bronze"
$ws.Range("D17").Value = 2014
$ws.Range("E17").Value = "Apr"
$ws.Range("F17").Value = "University of Nottingham – School of Medicine"
$ws.Range("G17").Value = "https://www.nottingham.ac.uk/medicine/about/athena-swan.aspx"
$ws.Range("H17").Value = "No"

# --- move the saved selection, matching the diff's sheetView/selection ---
$ws.Range("I19").Select() | Out-Null

# --- remove the now-empty "2014_Nov" sheet ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("2014_Nov").Delete() | Out-Null

# --- rename the remaining sheet ---
$ws.Name = "Selected"
